# Refactor sum/rest logic and generalize step handling
# Updates the "Paso 12.2" sheet so that the diagonal "carry" cells use
# generalized formulas (get_sum_num / get_rest_num style fractions) instead
# of hardcoded decimal literals, matching the refactor described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Paso 12.2")

# --- Row 3 (single-cell formula) -----------------------------------------
$ws.Range("C3").Formula = "=0.36/2"
$ws.Range("D3").Value = 0.01
$ws.Range("E3").Value = 0.05
$ws.Range("F3").Value = 0.01
$ws.Range("G3").Value = 0.05
$ws.Range("H3").Value = 0.01
$ws.Range("I3").Value = 0.01
$ws.Range("J3").Value = 0.01
$ws.Range("K3").Value = 0.01

# --- Row 4 (single-cell formula) -----------------------------------------
$ws.Range("C4").Value = 0.05
$ws.Range("D4").Formula = "=0.73/3"
$ws.Range("E4").Value = 0.01
$ws.Range("F4").Value = 0.01
$ws.Range("G4").Value = 0.05
$ws.Range("H4").Value = 0.05
$ws.Range("I4").Value = 0.01
$ws.Range("J4").Value = 0.01
$ws.Range("K4").Value = 0.01

# --- Row 5 (single-cell formula) -----------------------------------------
$ws.Range("C5").Value = 0.05
$ws.Range("D5").Value = 0.05
$ws.Range("E5").Formula = "=0.78/4"
$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0.05
$ws.Range("H5").Value = 0.05
$ws.Range("I5").Value = 0.05
$ws.Range("J5").Value = 0.01
$ws.Range("K5").Value = 0.01

# --- Row 6 (single-cell formula) -----------------------------------------
$ws.Range("C6").Value = 0.05
$ws.Range("D6").Value = 0.05
$ws.Range("E6").Value = 0.05
$ws.Range("F6").Formula = "=0.91/5"
$ws.Range("G6").Value = 0.05
$ws.Range("H6").Value = 0.05
$ws.Range("I6").Value = 0.05
$ws.Range("J6").Value = 0.05
$ws.Range("K6").Value = 0.01

# --- Row 7 -----------------------------------------------------------------
$ws.Range("C7").Value = 0.05
$ws.Range("D7").Value = 0.05
$ws.Range("E7").Value = 0.05
$ws.Range("F7").Value = 0.05
# G7:G11 becomes a shared formula group
$ws.Range("G7:G11").Formula = "=0.79/5"
$ws.Range("H7").Value = 0.01
$ws.Range("I7").Value = 0.01
$ws.Range("J7").Value = 0.01
$ws.Range("K7").Value = 0.01

# --- Row 8 -------------------------------------------------------------
$ws.Range("C8").Value = 0.18
# D8:D9 becomes a shared formula group
$ws.Range("D8:D9").Formula = "=0.73/3"
# E8:E10 becomes a shared formula group
$ws.Range("E8:E10").Formula = "=0.78/4"
# F8:F11 becomes a shared formula group
$ws.Range("F8:F11").Formula = "=0.91/5"
# G8 continues the G7:G11 shared group (set above)
# H8:H11 becomes a shared formula group
$ws.Range("H8:H11").Formula = "=0.82/5"
$ws.Range("I8").Value = 0.01
$ws.Range("J8").Value = 0.01
$ws.Range("K8").Value = 0.01

# --- Row 9 -------------------------------------------------------------
$ws.Range("C9").Value = 0.05
# D9, E9, F9, G9, H9 continue their shared groups set above
# I9:I11 becomes a shared formula group
$ws.Range("I9:I11").Formula = "=0.85/3"
$ws.Range("J9").Value = 0.01
$ws.Range("K9").Value = 0.01

# --- Row 10 ------------------------------------------------------------
$ws.Range("C10").Value = 0.05
$ws.Range("D10").Value = 0.05
# E10, F10, G10, H10, I10 continue their shared groups set above
# J10:J11 becomes a shared formula group
$ws.Range("J10:J11").Formula = "=0.88/2"
$ws.Range("K10").Value = 0.01

# --- Row 11 ------------------------------------------------------------
$ws.Range("C11").Value = 0.05
$ws.Range("D11").Value = 0.05
$ws.Range("E11").Value = 0.05
# F11, G11, H11, I11, J11 continue their shared groups set above
$ws.Range("K11").Value = 0.91

# Move the active selection on this sheet to mirror the author's last
# position after finishing the edits.
$ws.Range("M23").Select()
